$d = $word.ActiveDocument

# Merge the three runs (" avant que la méthode " / "OnDocumentChanged(" / ") ai fini son exécution")
# that were split around proofErr markers into a single run of continuous text,
# and remove the stray proofErr + bookmark markers.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(" avant que la méthode OnDocumentChanged() ai fini son exécution", $false, $false, $false, $false, $false, $true, 1, $false, " avant que la méthode OnDocumentChanged() ai fini son exécution", 2)
